# Actualización automática 2025-08-05 12:50:09
#
# A new asesor/client row ("ARCE CANDO DENISSE YAJAIRA") is inserted as the
# third data row (worksheet row 4) in both the "VENTAS POR GRUPO" and
# "VENTA MENSUAL" sheets, pushing the existing rows (and the trailing
# summary row) down by one. The new row gets all-zero figures, matching
# the pattern already used for rows with no sales. The summary/totals row
# at the bottom is updated to reflect the new total row count ("de 18" ->
# "de 19") on the sheet that tracks it textually.

$wb = $excel.ActiveWorkbook

# ---- Sheet "VENTAS POR GRUPO" --------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row before row 4 (shifts rows 4-20 down to 5-21).
$ws1.Rows.Item(4).Insert()

$ws1.Range("A4").Value = "OFICINA-CATAECSA"
$ws1.Range("B4").Value = "ARCE CANDO DENISSE YAJAIRA"
for ($c = 3; $c -le 18; $c++) {
  $ws1.Cells.Item(4, $c).Value = 0
}

# The old summary row (now row 21) counted "de 18" advisors; bump to 19.
for ($c = 3; $c -le 18; $c++) {
  $old = $ws1.Cells.Item(21, $c).Value()
  $prefix = $old.Substring(0, $old.IndexOf(" de "))
  $ws1.Cells.Item(21, $c).Value = "$prefix de 19"
}

# ---- Sheet "VENTA MENSUAL" ------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Insert a new row before row 4 (shifts rows 4-20 down to 5-21).
$ws2.Rows.Item(4).Insert()

$ws2.Range("A4").Value = "OFICINA-CATAECSA"
$ws2.Range("B4").Value = "ARCE CANDO DENISSE YAJAIRA"
for ($c = 3; $c -le 7; $c++) {
  $ws2.Cells.Item(4, $c).Value = 0
}
